$wb = $excel.ActiveWorkbook

# --- Sheet "RequestPayments" (sheet1.xml) ---
$wsRequest = $wb.Worksheets.Item("RequestPayments")
# Row 4 originally had A4/B4 raw values plus shared-formula anchors in C4/D4;
# rows 5-13 were the `#DIV/0!` shared-formula tail. Keep D4:D12 as empty,
# styled cells; drop A/B/C entirely and remove row 13.
$wsRequest.Range("A4:B4").Clear()
$wsRequest.Range("C4:C13").Clear()
$wsRequest.Range("D4:D13").ClearContents()
$wsRequest.Rows.Item(13).Delete()
$null = $wsRequest.Range("C4").Select()

# --- Sheet "Read payment details" (sheet3.xml) ---
$wsRead = $wb.Worksheets.Item("Read payment details")
# Remove stray rows 19, 13 and 11 (descending order so row numbers stay stable).
$wsRead.Rows.Item(19).Delete()
$wsRead.Rows.Item(13).Delete()
$wsRead.Rows.Item(11).Delete()

# --- Sheet "GetPaymentDetails" (sheet2.xml) ---
$wsGetPayment = $wb.Worksheets.Item("GetPaymentDetails")
$wsGetPayment.Visible = $false

# --- Sheet "Read payment details - 2" (sheet4.xml) ---
$wsRead2 = $wb.Worksheets.Item("Read payment details - 2")
$wsRead2.Visible = $false

# Select C13 on "Read payment details" and make it the active tab (also
# clears tabSelected from whichever sheet had it before, e.g. "Read payment
# details - 2").
$null = $wsRead.Range("C13").Select()
$wsRead.Activate()
